$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AY2").Value = 28898619.392
$ws.Range("AZ2").Value = 28512851.968
$ws.Range("AY3").Value = 7843675.136
$ws.Range("AZ3").Value = 6967942.144
$ws.Range("AY4").Value = 1133246.976
$ws.Range("AZ4").Value = 1601468.032
$ws.Range("AY5").Value = 2035003.008
$ws.Range("AZ5").Value = 2425200.896
$ws.Range("AY6").Value = 1640957.056
$ws.Range("AZ6").Value = 2025688.96
$ws.Range("AY7").Value = 425951.008
$ws.Range("AZ7").Value = 366280
$ws.Range("AY9").Value = 275577.984
$ws.Range("AZ9").Value = 213858
$ws.Range("AY11").Value = 2332939.008
$ws.Range("AZ11").Value = 335446.016
$ws.Range("AY12").Value = 2041001.984
$ws.Range("AZ12").Value = 1766494.976
$ws.Range("AY15").Value = 13817
$ws.Range("AZ15").Value = 14237
$ws.Range("AY16").Value = 335862.016
$ws.Range("AZ16").Value = 251587.008
$ws.Range("AY19").Value = 915507.008
$ws.Range("AZ19").Value = 904160
$ws.Range("AY21").Value = 281129.984
$ws.Range("AY22").Value = 953
$ws.Range("AZ22").Value = 1211
$ws.Range("AY23").Value = 4146097.92
$ws.Range("AZ23").Value = 4201251.072
$ws.Range("AY24").Value = 14866892.8
$ws.Range("AZ24").Value = 15575954.432
$ws.Range("AY26").Value = 28898619.392
$ws.Range("AZ26").Value = 28512851.968
$ws.Range("AY27").Value = 5095707.136
$ws.Range("AZ27").Value = 4087998.976
$ws.Range("AY28").Value = 412902.016
$ws.Range("AZ28").Value = 387081.984
$ws.Range("AY29").Value = 710553.9840000001
$ws.Range("AZ29").Value = 964220.992
$ws.Range("AY30").Value = 136514
$ws.Range("AZ30").Value = 143668.992
$ws.Range("AY31").Value = 2126443.008
$ws.Range("AZ31").Value = 2258499.072
$ws.Range("AY32").Value = 0
$ws.Range("AY34").Value = 290423.008
$ws.Range("AZ34").Value = 334528
$ws.Range("AY36").Value = 1418871.04
$ws.Range("AY37").Value = 10077039.616
$ws.Range("AZ37").Value = 10654189.568
$ws.Range("AY38").Value = 7088996.864
$ws.Range("AZ38").Value = 7635253.76
$ws.Range("AY40").Value = 184656
$ws.Range("AZ40").Value = 271103.008
$ws.Range("AY41").Value = 485479.008
$ws.Range("AZ41").Value = 669257.9840000001
$ws.Range("AY43").Value = 2317907.968
$ws.Range("AZ43").Value = 2078574.976
$ws.Range("AY46").Value = 1046888
$ws.Range("AZ46").Value = 1044073.984
$ws.Range("AY47").Value = 12678986.176
$ws.Range("AZ47").Value = 12726589.952
$ws.Range("AY48").Value = 7667615.232
$ws.Range("AZ48").Value = 7667615.232
$ws.Range("AY49").Value = 5565511.168
$ws.Range("AZ49").Value = 5058975.232
$ws.Range("AY51").Value = 0
$ws.Range("AZ51").Value = 0
$ws.Range("AY52").Value = -554140.032
$ws.Range("AZ52").Value = 0
$ws.Range("AY59").Value = 1497580.288
$ws.Range("AZ59").Value = 1176647.936
$ws.Range("AY60").Value = -621199.936
$ws.Range("AZ60").Value = -431537.984
$ws.Range("AY61").Value = 876379.968
$ws.Range("AZ61").Value = 745110.0159999999
$ws.Range("AY62").Value = -90452.024
$ws.Range("AZ62").Value = -123488
$ws.Range("AY63").Value = -503431.04
$ws.Range("AZ63").Value = -362376.992
$ws.Range("AY64").Value = 352606.016
$ws.Range("AY65").Value = 3410
$ws.Range("AZ65").Value = 2096
$ws.Range("AY66").Value = -537654.0159999999
$ws.Range("AZ66").Value = -102140
$ws.Range("AY67").Value = 258
$ws.Range("AZ67").Value = 574
$ws.Range("AY68").Value = -155336.032
$ws.Range("AZ68").Value = -191120.992
$ws.Range("AY69").Value = 92094.984
$ws.Range("AZ69").Value = 150564.992
$ws.Range("AY70").Value = -247431.024
$ws.Range("AZ70").Value = -341686.016
$ws.Range("AY74").Value = -54219
$ws.Range("AZ74").Value = -31346
$ws.Range("AY75").Value = 4864
$ws.Range("AZ75").Value = 10976
$ws.Range("AY76").Value = -53959
$ws.Range("AZ76").Value = 11781
$ws.Range("AY79").Value = -5318
$ws.Range("AZ79").Value = -4518
$ws.Range("AY80").Value = 65015.024
$ws.Range("AZ80").Value = -13107
